# Fruta / hortaliza, semanal
# A new weekly price record was added to the top of the data (row 13),
# pushing the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; everything from row 13 down (through
# the old last row, 108) shifts down to make room (new last row: 109).
$ws.Rows("13:13").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(13, 1).Value = 3
$ws.Cells.Item(13, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = 44537
$ws.Cells.Item(13, 5).Value = 5
$ws.Cells.Item(13, 6).Value = 100112026
$ws.Cells.Item(13, 7).Value = "Haba"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 70
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 8000
$ws.Cells.Item(13, 13).Value = 8000
$ws.Cells.Item(13, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(13, 16).Value = 320
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(13, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
